$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = "LM-White-Paper-Threat-Driven-Approach.pdf"
$ws.Range("B11").Value = "pdf"
$ws.Range("C11").Value = "uconn-sdp-team11-tagged-docs"
$ws.Range("D11").Value = "defense"
$ws.Range("E11").Value = "lm"
$ws.Range("F11").Value = "threat"
